$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.3044033944606781
$ws.Range("B1").Value = 1.807124733924866
$ws.Range("C1").Value = 3.61018443107605
$ws.Range("D1").Value = 3.397135257720947
$ws.Range("E1").Value = 0.829828143119812
